$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target item names (column A) and quantities (column B) for rows 2..110
$names = @(
  "ВердиоГаст® Растительный комплекс для улучшения пищеварения (БАД ),  капсулы"
  "Полынь горькая трава 50г"
  "Шиповник плоды низковитаминные 50г"
  "Липа цветки 35г"
  "Кукуруза столбики с рыльцами 40г"
  "Сб. Фитонефрол (Урологический сбор) 50г"
  "Укроп пахучий плоды 50г"
  "Брусника листья 50г"
  "Эвкалипт прутовидный листья 75г"
  "Дуба кора 75г"
  "Ромашка цветки вн 50г"
  "Мята перечная листья 50г"
  "Ноготки цветки 50г"
  "Эрва шерстистая трава 30г"
  "Береза почки 50г"
  "Сб. Грудной №4 50г"
  "Багульник болотный побеги 50г"
  "Чабрец трава 50г"
  "Боярышник плоды 75г"
  "Валериана корневища с корнями 50г"
  "Девясил корневища и корни 50г"
  "Пустырник трава 50г"
  "Бессмертник песчаный цветки 30г"
  "Шалфей листья 50г"
  "Подорожник большой листья 50г"
  "Ламинарии слоевища (морская капуста) 100г"
  "Сб. Фитогепатол №2 (Желчегонный сбор №2) 35г"
  "Мать-и-мачеха листья 35г"
  "Аир корневища 75г"
  "Пижма цветки 75г"
  "Лен семена 100г"
  "Рябина плоды 50г"
  "Череда трава 50г"
  "Крушина кора 50г"
  "Зверобой трава 50г"
  "Можжевельник плоды 50г"
  "Тысячелистник трава 50г"
  "Солодка корни 50г"
  "Толокнянка листья 50г"
  "Сенна листья 50г"
  "Сб. Фитопектол №1 (Грудной сбор №1) 35г"
  "Спорыш трава 50г"
  "Алтей корни 75г"
  "Чага (березовый гриб) 50г"
  "Крапива листья 50г"
  "Чистотел трава 50г"
  "Сб. Фитопектол №2 (Грудной сбор №2) 35г"
  "Фп Детский травяной чай ""ФармаЦветик® для иммунитета"" 20х1,5 г"
  "Фп Фиточай ""Лактафитол"" (БАД) 20х1,5 г"
  "Фп Детский травяной чай ""ФармаЦветик®  при простуде"" 20х1,5 г"
  "Фп Детский травяной чай ""ФармаЦветик® для животика"" 20х1,5 г"
  "Фп Детский травяной чай ""ФармаЦветик® для спокойного сна"" 20х1,5 г"
  "Фп ""ВердиоГаст® Фиточай для улучшения пищеварения с зеленым чаем""(БАД) 20*1,5г"
  "Фп ""ВердиоГаст® Фиточай для улучшения пищеварения с черным чаем"" (БАД) 20*1,5г"
  "Фп ""Щедрость природы® Фиточай для иммунитета"" 20х2,0 г"
  "Фп ""Щедрость природы® Фиточай при простуде"" 20х2,0 г"
  "Фп ""Щедрость природы® Фиточай кардиологический"" 20х2,0 г"
  "Фп ""Щедрость природы® Фиточай успокоительный""20х2,0 г"
  "Фп Шалфей листья 20х1,5г"
  "Фп ""Щедрость природы® Фиточай диабетический"" 20х2,0 г"
  "Фп Сб. Фитоседан №2 (Успокоительный сбор №2) 20x2,0г"
  "Фп Сб. Бруснивер 20x2,0г"
  "Фп Сб. Грудной №4 20x2,0г"
  "Фп Сб. Проктофитол (Противогеморроидальный сбор) 20х2,0г"
  "Фп Толокнянка листья 20x1,5г"
  "Фп Сб. Фитогастрол (Желудочно-кишечный сбор) 20x2,0г"
  "Фп Мелисса лекарственная трава 20x1,5г"
  "Фп Сенна листья 20x1,5г"
  "Фп Липа цветки 20x1,5г"
  "Фп Аир корневища 20x1,5г"
  "Фп Боярышник плоды 20х3,0г"
  "Фп Сб. Фитоседан №3 (Успокоительный сбор №3) 20х2,0г"
  "Фп Пижма цветки 20х1,5г"
  "Фп Сб. Фитогепатол №3 (Желчегонный сбор №3) 20x2,0г"
  "Фп Сб. Элекасол 20x2,0г"
  "Фп Фиточай ""Тибетский"" (БАД) 20х2,0  г"
  "Фп Фиточай ""Опалиховский"" (БАД) 20х2,0 г"
  "Фп Ромашка цветки 20x1,5г"
  "Фп ""Щедрость природы® Фиточай очищающий"" 20х2,0 г"
  "Фп Сб. Арфазетин-Э 20x2,0г"
  "Фп Мята перечная листья 20x1,5г"
  "Фп Чистотел трава 20х1,5г"
  "Фп Шиповник плоды 20х2,0г"
  "Фп Пустырник трава 20x1,5г"
  "Фп Подорожник листья 20x1,5г"
  "Фп Брусника листья 20х1,5г"
  "Фп Зверобой трава 20x1,5г"
  "Фп ""Щедрость природы® Фиточай для пищеварения"" 20х2,0 г"
  "Фп Череда трава 20х1,5г"
  "Фп Сб. Фитонефрол (Урологический сбор) 20x2,0г"
  "Фп Чабрец трава 20x1,5 г"
  "Фп Душица трава 20x1,5г"
  "Фп Крапива листья 20x1,5г"
  "Фп Пастушья сумка трава 20х1,5г"
  "Фп Хвощ полевой трава 20х1,5г"
  "Фп Сб. Желудочный №3 20x2,0г"
  "Фп Береза листья 20x1,5г"
  "Фп Золототысячник трава 20х1,5г"
  "Фп Фиалка трехцветная трава 20x1,5г"
  "Фп Ольха соплодия 20х1,5г"
  "Фп Валериана корневища с корнями 20x1,5г"
  "Фп Тысячелистник трава 20x1,5г"
  "Фп Лапчатка корневища 20x2,5г"
  "Фп Девясил корневища и корни 20х1,5г"
  "Фп Крушина кора 20x1,5г"
  "Фп Ноготки цветки 20x1,5г"
  "Фп Почечный чай листья 20x1,5г"
  "Фп Бадан корневища 20x1,5г"
  "Фп Кровохлебка корневища и корни 20x1,5г"
)

$values = @(
  81472
  24836
  21285
  14118
  19252
  6317
  48692
  12359
  21733
  52997
  84567
  19567
  21531
  12388
  17031
  33470
  13277
  20329
  21256
  20158
  18419
  12068
  29235
  37990
  10010
  17944
  4943
  31754
  8303
  20097
  65726
  2128
  14311
  12536
  46695
  18380
  20815
  51436
  11227
  34733
  8539
  22647
  9737
  45668
  23133
  29771
  11292
  1208
  10529
  3271
  3900
  5888
  7990
  9400
  1008
  1026
  1476
  1962
  96810
  1062
  37305
  139827
  548697
  17380
  26766
  67275
  34399
  60940
  62296
  3406
  15956
  80763
  5664
  84274
  41666
  7200
  4932
  1311041
  1800
  42641
  73126
  33486
  55991
  46365
  32545
  85913
  61495
  1854
  57564
  220855
  87279
  35946
  84526
  6640
  36754
  31431
  6610
  6433
  6088
  6181
  15279
  16684
  3466
  13614
  9513
  43788
  73644
  2041
  12454
)

for ($i = 0; $i -lt $names.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 1).Value = $names[$i]
  $ws.Cells.Item($row, 2).Value = $values[$i]
}

# Update the view: scroll position and active-cell selection
$ws.Activate()
$ws.Range("A94").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 82
$win.ScrollColumn = 1
